$d = $word.ActiveDocument

$replacements = @(
    @{old="168÷7="; new="127÷2="},
    @{old="195÷8="; new="576÷2="},
    @{old="140÷3="; new="501÷8="},
    @{old="658÷6="; new="764÷6="},
    @{old="849÷5="; new="326÷5="},
    @{old="803÷9="; new="638÷4="},
    @{old="616÷2="; new="126÷3="},
    @{old="961÷5="; new="892÷8="},
    @{old="554÷9="; new="107÷8="},
    @{old="330÷5="; new="759÷5="},
    @{old="946÷7="; new="245÷7="},
    @{old="592÷6="; new="127÷5="},
    @{old="156÷7="; new="547÷4="},
    @{old="338÷6="; new="819÷2="},
    @{old="894÷9="; new="499÷3="},
    @{old="142÷8="; new="702÷4="},
    @{old="395÷9="; new="807÷5="},
    @{old="552÷9="; new="509÷8="},
    @{old="832÷5="; new="955÷3="},
    @{old="127÷6="; new="815÷3="},
    @{old="644÷4="; new="197÷9="},
    @{old="788÷5="; new="307÷6="},
    @{old="887÷2="; new="572÷6="},
    @{old="648÷3="; new="646÷3="},
    @{old="675÷6="; new="519÷4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
